# Generate Report for Handback
# Fill in the "Latest Target File / Latest Handback File / Latest Handback
# DateTime / Error Detail" columns for the row that tracks
# 8cc0f32d-1998-46a8-846f-a4a610794971.md on both the zh-cn and de-de
# report sheets, since a fresh handback was generated for that file but its
# content is stale compared to the latest handoff.

$wb = $excel.ActiveWorkbook

$targetMd = "8cc0f32d-1998-46a8-846f-a4a610794971.md"
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7fecc50e3d87a280e3ea537bf085d6ddea03729e/e2e/8cc0f32d-1998-46a8-846f-a4a610794971.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2994da280c7fc17c17f555d81ebad97ab8a9146c/e2e/8cc0f32d-1998-46a8-846f-a4a610794971.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7fecc50e3d87a280e3ea537bf085d6ddea03729e/e2e/8cc0f32d-1998-46a8-846f-a4a610794971.md."

# zh-cn sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("J7").Value = $wsZh.Range("G7").Value2
$wsZh.Range("K7").Value = "2016-09-07 07:18:18"
$wsZh.Range("P7").Value = $errorDetail
$wsZh.Range("I7").Value = $targetMd
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $targetUrl, "", "", $targetMd)
$wsZh.Range("I7").Style = "HyperLink"

# de-de sheet
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("J7").Value = $wsDe.Range("G7").Value2
$wsDe.Range("K7").Value = "2016-09-07 07:18:27"
$wsDe.Range("P7").Value = $errorDetail
$wsDe.Range("I7").Value = $targetMd
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $targetUrl, "", "", $targetMd)
$wsDe.Range("I7").Style = "HyperLink"
